# Update for 2021 preliminary elections
# Append new file-index rows to the bottom of Sheet1's data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newEntries = @(
    "data-raw/elections/Official Results September 14, 2021.pdf",
    "data-raw/elections/Official Results September 14, 2021.xlsx",
    "data/elections/2021_preliminary_turnout_bydistrict_framingham.csv",
    "data/elections/tidy/2021_preliminary_turnout_tidy.csv",
    "data/elections/2021_preliminary_turnout.csv"
)

# Find the first empty row in column A below the existing data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newEntries.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newEntries[$i]
}

# Keep the active selection / view consistent with the new last row.
$lastDataRow = $startRow + $newEntries.Count - 1
$ws.Range("A" + $lastDataRow).Select() | Out-Null
